$d = $word.ActiveDocument

# Step 1: temporarily mark the "enviados para" run as bold. Giving it a
# distinct run identity (different formatting) from its neighbours keeps
# the engine from folding it into whichever run ends up adjacent to it
# once the intervening runs are edited/removed below.
$rngMark = $d.Content
$rngMark.Find.Execute(
    "enviados para", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
$rngMark.Bold = 1

# Step 2: shorten the lead-in text of the first run.
$rngText = $d.Content
$rngText.Find.Execute(
    "Todos os trabalhos deverão ser entregues através do sistema ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Todos os trabalhos deverão ser ", 2
)

# Step 3: delete the runs "Microsoft Teams", " ", "(Chat Geral) " and "ou "
# (they sit between "...ser " and "enviados para").
$rngDel = $d.Content
$rngDel.Find.Execute(
    "Microsoft Teams (Chat Geral) ou ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
$rngDel.Delete()

# Step 4: restore "enviados para" to its normal (non-bold) formatting.
$rngUnmark = $d.Content
$rngUnmark.Find.Execute(
    "enviados para", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
$rngUnmark.Bold = 0
